$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (本次监测期内均值/上周期均值 offset) measurement values, and
# column D where the underlying source re-measurement also changed.
# Cells that previously held the literal "--" placeholder are switched to a
# real numeric reading (or cleared) where the new data now provides one.

$ws.Range("G2").Value = -0.003
$ws.Range("G3").Value = -0.014
$ws.Range("G4").Value = -0.003

$ws.Range("D5").Value = ""
$ws.Range("G5").Value = 0.119

$ws.Range("G6").Value = 0.004
$ws.Range("G7").Value = 0.017

$ws.Range("D8").Value = ""
$ws.Range("G8").Value = ""

$ws.Range("G9").Value = -0.002
$ws.Range("G10").Value = -0.002
$ws.Range("G11").Value = 0.014
$ws.Range("G12").Value = -0.029
$ws.Range("G13").Value = -0.021

$ws.Range("D14").Value = ""
$ws.Range("G14").Value = 0.646

$ws.Range("G15").Value = 0.008

$ws.Range("D16").Value = ""
$ws.Range("G16").Value = -0.01

$ws.Range("D17").Value = 0.071
$ws.Range("G17").Value = 0.061

$ws.Range("D18").Value = ""
$ws.Range("G18").Value = -0.061

$ws.Range("G19").Value = -0.006

$ws.Range("D20").Value = ""
$ws.Range("G20").Value = -0.012

$ws.Range("G21").Value = -0.015
$ws.Range("G22").Value = -0.014

$ws.Range("D23").Value = ""
$ws.Range("G23").Value = 0.018

$ws.Range("G24").Value = -0.017
$ws.Range("G25").Value = -0.01
$ws.Range("G26").Value = 0.003
$ws.Range("G27").Value = -0.02
$ws.Range("G29").Value = -0.023
$ws.Range("G30").Value = -0.006
$ws.Range("G31").Value = 0.02

$ws.Range("D32").Value = ""
$ws.Range("G32").Value = -0.008

$ws.Range("G33").Value = 0.009
$ws.Range("G34").Value = 0.017
$ws.Range("G35").Value = 0.015

$ws.Range("D37").Value = ""
$ws.Range("G37").Value = 0.181

$ws.Range("G38").Value = 0.004
$ws.Range("G39").Value = -0.022

$ws.Range("D40").Value = ""
$ws.Range("G40").Value = -0.343

$ws.Range("G41").Value = -0.015
$ws.Range("G42").Value = -0.001

$ws.Range("D43").Value = 0.026
$ws.Range("G43").Value = 0.015

$ws.Range("G45").Value = -0.032
$ws.Range("G46").Value = -0.005
$ws.Range("G47").Value = -0.004
$ws.Range("G48").Value = -0.001
$ws.Range("G49").Value = -0.005

# Reflect the user's final selection: column G for all data rows.
$ws.Range("G2:G49").Select()
